$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-7 from 45243 (2023-11-13)
# to 45244 (2023-11-14).
$ws.Range("C2:C7").Value = 45244
